# Update the loading_percent results (rows 2-25, cols B-O) for the
# "case with 380 kV done" run. Columns D and I stay 0 (unchanged);
# all other columns get new simulation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 24,14
$values[0,0] = 12.38884992330517
$values[0,1] = 12.54581670815302
$values[0,2] = 0
$values[0,3] = 16.99707910483048
$values[0,4] = 38.35872261492172
$values[0,5] = 32.65341674856203
$values[0,6] = 15.41191880399084
$values[0,7] = 0
$values[0,8] = 8.159328671494469
$values[0,9] = 8.002756177034765
$values[0,10] = 12.35441517174404
$values[0,11] = 15.61439809887665
$values[0,12] = 20.15170681335833
$values[0,13] = 23.94871825159753
$values[1,0] = 12.182900690301
$values[1,1] = 12.57499759424271
$values[1,2] = 0
$values[1,3] = 17.01602150037447
$values[1,4] = 38.40170507286052
$values[1,5] = 32.75646789101586
$values[1,6] = 15.45601497448684
$values[1,7] = 0
$values[1,8] = 8.14357779226928
$values[1,9] = 7.845717438373098
$values[1,10] = 12.35144222344896
$values[1,11] = 15.5769369681318
$values[1,12] = 20.20606974307067
$values[1,13] = 24.02524754820314
$values[2,0] = 12.05665463757624
$values[2,1] = 12.59386898834525
$values[2,2] = 0
$values[2,3] = 17.02992135315145
$values[2,4] = 38.43575774166634
$values[2,5] = 32.82742948173702
$values[2,6] = 15.48500033076746
$values[2,7] = 0
$values[2,8] = 8.133719458811397
$values[2,9] = 7.748654586277483
$values[2,10] = 12.35112123201859
$values[2,11] = 15.55588709661162
$values[2,12] = 20.24105889741367
$values[2,13] = 24.07613855620889
$values[3,0] = 12.0053257985865
$values[3,1] = 12.60179983510527
$values[3,2] = 0
$values[3,3] = 17.03615700034855
$values[3,4] = 38.45156085998552
$values[3,5] = 32.8582755670365
$values[3,6] = 15.49729292759783
$values[3,7] = 0
$values[3,8] = 8.129655400281605
$values[3,9] = 7.708991118863493
$values[3,10] = 12.35137023258481
$values[3,11] = 15.54780578179775
$values[3,12] = 20.25572325403168
$values[3,13] = 24.09785798008969
$values[4,0] = 11.99681169682676
$values[4,1] = 12.60313129717007
$values[4,2] = 0
$values[4,3] = 17.03722695284632
$values[4,4] = 38.4543012894117
$values[4,5] = 32.86351388626736
$values[4,6] = 15.49936315995382
$values[4,7] = 0
$values[4,8] = 8.128977756225419
$values[4,9] = 7.702399943790285
$values[4,10] = 12.35143456017437
$values[4,11] = 15.54649405155186
$values[4,12] = 20.2581828132085
$values[4,13] = 24.10152371363733
$values[5,0] = 12.05596183666964
$values[5,1] = 12.59397497145909
$values[5,2] = 0
$values[5,3] = 17.03000313498533
$values[5,4] = 38.43596306879199
$values[5,5] = 32.82783767987604
$values[5,6] = 15.48516416548678
$values[5,7] = 0
$values[5,8] = 8.133664838127444
$values[5,9] = 7.748120045497593
$values[5,10] = 12.35112305039189
$values[5,11] = 15.55577609050469
$values[5,12] = 20.24125502071922
$values[5,13] = 24.07642750048461
$values[6,0] = 12.31783445202839
$values[6,1] = 12.55568062667336
$values[6,2] = 0
$values[6,3] = 17.00313995443027
$values[6,4] = 38.37195297536663
$values[6,5] = 32.68735055781622
$values[6,6] = 15.4267270238416
$values[6,7] = 0
$values[6,8] = 8.153936850348362
$values[6,9] = 7.948772572916868
$values[6,10] = 12.35307879762047
$values[6,11] = 15.60108017395888
$values[6,12] = 20.17011765327837
$values[6,13] = 23.97429564240671
$values[7,0] = 12.83019295009828
$values[7,1] = 12.48812672569453
$values[7,2] = 0
$values[7,3] = 16.96843090622716
$values[7,4] = 38.30720191415585
$values[7,5] = 32.47304612950107
$values[7,6] = 15.32726497671237
$values[7,7] = 0
$values[7,8] = 8.192189184493207
$values[7,9] = 8.334988342126294
$values[7,10] = 12.36878086197829
$values[7,11] = 15.70512869171441
$values[7,12] = 20.04334097065694
$values[7,13] = 23.80498194400749
$values[8,0] = 13.20222899821678
$values[8,1] = 12.4430499601404
$values[8,2] = 0
$values[8,3] = 16.95383446329299
$values[8,4] = 38.29662447956206
$values[8,5] = 32.35314710446899
$values[8,6] = 15.26338395167545
$values[8,7] = 0
$values[8,8] = 8.219358801098153
$values[8,9] = 8.61152559295733
$values[8,10] = 12.387452299428
$values[8,11] = 15.79045449268326
$values[8,12] = 19.95788094936479
$values[8,13] = 23.69947592641854
$values[9,0] = 13.36978225275151
$values[9,1] = 12.4235235007072
$values[9,2] = 0
$values[9,3] = 16.94954935654398
$values[9,4] = 38.29982099481899
$values[9,5] = 32.30680267872511
$values[9,6] = 15.23631255107111
$values[9,7] = 0
$values[9,8] = 8.231511926482531
$values[9,9] = 8.735216742911629
$values[9,10] = 12.39747213062963
$values[9,11] = 15.83111019582092
$values[9,12] = 19.92065551821522
$values[9,13] = 23.6555822370687
$values[10,0] = 13.43292895844154
$values[10,1] = 12.41626947497909
$values[10,2] = 0
$values[10,3] = 16.94826407779072
$values[10,4] = 38.30217966989543
$values[10,5] = 32.29043548872561
$values[10,6] = 15.22634677291281
$values[10,7] = 0
$values[10,8] = 8.236083903239038
$values[10,9] = 8.781710170011129
$values[10,10] = 12.40148359076144
$values[10,11] = 15.84676199934775
$values[10,12] = 19.90679543941126
$values[10,13] = 23.63955086417895
$values[11,0] = 13.41934353637541
$values[11,1] = 12.41782553211386
$values[11,2] = 0
$values[11,3] = 16.94852590042867
$values[11,4] = 38.3016206804099
$values[11,5] = 32.29390781409639
$values[11,6] = 15.22848038616205
$values[11,7] = 0
$values[11,8] = 8.235100594438672
$values[11,9] = 8.771713007375839
$values[11,10] = 12.40061003019473
$values[11,11] = 15.84337983648829
$values[11,12] = 19.90976995782439
$values[11,13] = 23.64297725069815
$values[12,0] = 13.37498375454761
$values[12,1] = 12.42292390074369
$values[12,2] = 0
$values[12,3] = 16.94943686135462
$values[12,4] = 38.29999205405002
$values[12,5] = 32.30543242331313
$values[12,6] = 15.23548693919053
$values[12,7] = 0
$values[12,8] = 8.231888672618471
$values[12,9] = 8.739048949459608
$values[12,10] = 12.39779781704567
$values[12,11] = 15.83239278963491
$values[12,12] = 19.91951050935151
$values[12,13] = 23.65425149524282
$values[13,0] = 13.34777105303781
$values[13,1] = 12.42606504234902
$values[13,2] = 0
$values[13,3] = 16.95003875207456
$values[13,4] = 38.29914389440219
$values[13,5] = 32.31264566661207
$values[13,6] = 15.23981583493916
$values[13,7] = 0
$values[13,8] = 8.229917337900069
$values[13,9] = 8.718995023485956
$values[13,10] = 12.39610346694061
$values[13,11] = 15.82569604195703
$values[13,12] = 19.92550763330857
$values[13,13] = 23.66123417112835
$values[14,0] = 13.19123989868508
$values[14,1] = 12.44434571376887
$values[14,2] = 0
$values[14,3] = 16.9541617885602
$values[14,4] = 38.29657645589549
$values[14,5] = 32.35634110203488
$values[14,6] = 15.26519311245922
$values[14,7] = 0
$values[14,8] = 8.218560372060937
$values[14,9] = 8.603395973276401
$values[14,10] = 12.38682799664319
$values[14,11] = 15.78783388099898
$values[14,12] = 19.96034685119874
$values[14,13] = 23.70242705172195
$values[15,0] = 13.09473771887669
$values[15,1] = 12.45581068855789
$values[15,2] = 0
$values[15,3] = 16.95729343805207
$values[15,4] = 38.29705036042183
$values[15,5] = 32.38524923384378
$values[15,6] = 15.2812702713088
$values[15,7] = 0
$values[15,8] = 8.211540035583271
$values[15,9] = 8.531908800578627
$values[15,10] = 12.38152708291958
$values[15,11] = 15.7650722739049
$values[15,12] = 19.98214167210305
$values[15,13] = 23.72874835299766
$values[16,0] = 13.03907673841531
$values[16,1] = 12.46249723917542
$values[16,2] = 0
$values[16,3] = 16.95931641356522
$values[16,4] = 38.29807695713819
$values[16,5] = 32.4026478941235
$values[16,6] = 15.29070460604733
$values[16,7] = 0
$values[16,8] = 8.207482717607215
$values[16,9] = 8.490595337666111
$values[16,10] = 12.37862192620703
$values[16,11] = 15.75215415848361
$values[16,12] = 19.99483291275699
$values[16,13] = 23.74427370027948
$values[17,0] = 13.02020603653152
$values[17,1] = 12.4647770437555
$values[17,2] = 0
$values[17,3] = 16.96003947708224
$values[17,4] = 38.2985541391457
$values[17,5] = 32.4086711681767
$values[17,6] = 15.29393106868372
$values[17,7] = 0
$values[17,8] = 8.206105660295325
$values[17,9] = 8.476575002621747
$values[17,10] = 12.37766305338243
$values[17,11] = 15.74781039789582
$values[17,12] = 19.99915667553786
$values[17,13] = 23.74959660425346
$values[18,0] = 13.10502707087496
$values[18,1] = 12.45458068429012
$values[18,2] = 0
$values[18,3] = 16.95693712862265
$values[18,4] = 38.29692189631814
$values[18,5] = 32.38209204309992
$values[18,6] = 15.27953946066042
$values[18,7] = 0
$values[18,8] = 8.212289374144333
$values[18,9] = 8.539539328282613
$values[18,10] = 12.38207650878742
$values[18,11] = 15.76747736163093
$values[18,12] = 19.97980549638397
$values[18,13] = 23.72590644969167
$values[19,0] = 13.38802195549882
$values[19,1] = 12.4214225853762
$values[19,2] = 0
$values[19,3] = 16.94916014349803
$values[19,4] = 38.30043928886882
$values[19,5] = 32.30201524873792
$values[19,6] = 15.23342119713131
$values[19,7] = 0
$values[19,8] = 8.232832913803694
$values[19,9] = 8.748652885528664
$values[19,10] = 12.39861795718953
$values[19,11] = 15.8356130613784
$values[19,12] = 19.91664306571885
$values[19,13] = 23.65092395480439
$values[20,0] = 13.57118918586542
$values[20,1] = 12.40056889710918
$values[20,2] = 0
$values[20,3] = 16.94604340301343
$values[20,4] = 38.30942877726458
$values[20,5] = 32.25657397757774
$values[20,6] = 15.20494462854609
$values[20,7] = 0
$values[20,8] = 8.246083362020634
$values[20,9] = 8.883287566339671
$values[20,10] = 12.4106934188393
$values[20,11] = 15.88163395570059
$values[20,12] = 19.87674014026219
$values[20,13] = 23.60535896573569
$values[21,0] = 13.47361172737602
$values[21,1] = 12.41162433500513
$values[21,2] = 0
$values[21,3] = 16.94752741395978
$values[21,4] = 38.30402001342941
$values[21,5] = 32.2801950024194
$values[21,6] = 15.2199909334894
$values[21,7] = 0
$values[21,8] = 8.239027599987939
$values[21,9] = 8.811629970331554
$values[21,10] = 12.40413358796203
$values[21,11] = 15.85693819424579
$values[21,12] = 19.89791137802008
$values[21,13] = 23.62936292599089
$values[22,0] = 13.10037581278885
$values[22,1] = 12.45513647269097
$values[22,2] = 0
$values[22,3] = 16.95709752281535
$values[22,4] = 38.29697762561909
$values[22,5] = 32.3835169828177
$values[22,6] = 15.28032136404721
$values[22,7] = 0
$values[22,8] = 8.211950664284004
$values[22,9] = 8.536090231656127
$values[22,10] = 12.3818276695952
$values[22,11] = 15.76638949805413
$values[22,12] = 19.98086117950789
$values[22,13] = 23.72719005041345
$values[23,0] = 12.69209870575229
$values[23,1] = 12.50559889700714
$values[23,2] = 0
$values[23,3] = 16.97590168093399
$values[23,4] = 38.31821555256327
$values[23,5] = 32.52444362110635
$values[23,6] = 15.35255531618741
$values[23,7] = 0
$values[23,8] = 8.182003576582506
$values[23,9] = 8.231580463888093
$values[23,10] = 12.36327234538236
$values[23,11] = 15.67539160562704
$values[23,12] = 20.07628295302818
$values[23,13] = 23.84746945933033

$ws.Range("B2:O25").Value = $values
